$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Refresh the cached "last rendered" date/time field text
#    (datetimeFigureOut / datetime1 fields) from 1/8/2024 to 1/17/2024
#    across the slide master, every slide layout, the handout master
#    and the notes master - mirroring PowerPoint's automatic refresh
#    of these cached field values when the deck is reopened/resaved.
# ------------------------------------------------------------------

$newDate = "1/17/2024"

$design = $p.Designs.Item(1)
$master = $design.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $sh = $master.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $cl = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $cl.Shapes.Count; $j++) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

$hm = $p.HandoutMaster
for ($j = 1; $j -le $hm.Shapes.Count; $j++) {
    $sh = $hm.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

$nm = $p.NotesMaster
for ($j = 1; $j -le $nm.Shapes.Count; $j++) {
    $sh = $nm.Shapes.Item($j)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# ------------------------------------------------------------------
# 2) Remove the stray "Picture 6" image that was pasted onto slide 2
#    (Multiple Inheritance slide).
# ------------------------------------------------------------------

$s2 = $p.Slides.Item(2)
for ($j = $s2.Shapes.Count; $j -ge 1; $j--) {
    $sh = $s2.Shapes.Item($j)
    if ($sh.Name -eq "Picture 6") {
        $sh.Delete()
    }
}
